$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "'70"
$ws.Range("E4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("L4").Value = "'65"
$ws.Range("M4").Value = 0
$ws.Range("AB4").Value = "'110"
$ws.Range("AC4").Value = 60
$ws.Range("D5").Value = "'8"
$ws.Range("E5").Value = "'8"
$ws.Range("L5").Value = "'12"
$ws.Range("M5").Value = "'12"
$ws.Range("AB5").Value = "'14"
$ws.Range("AC5").Value = "'14"
$ws.Range("D6").Value = 8
$ws.Range("L6").Value = 7
$ws.Range("AB6").Value = 77
$ws.Range("AC7").Value = 0.59
$ws.Range("D8").Value = "'1"
$ws.Range("L8").Value = "'1"
$ws.Range("AB8").Value = "'1"
$ws.Range("D9").Value = "'5"
$ws.Range("E9").Value = 0
$ws.Range("L9").Value = "'2"
$ws.Range("AB9").Value = "'8"
$ws.Range("E10").Value = 0
$ws.Range("D11").Value = "'8"
$ws.Range("L11").Value = "'8"
$ws.Range("AB11").Value = "'12"
$ws.Range("D12").Value = "'8"
$ws.Range("L12").Value = "'6"
$ws.Range("AB12").Value = "'10"
$ws.Range("D13").Value = "'8"
$ws.Range("L13").Value = "'6"
$ws.Range("AB13").Value = "'7"
$ws.Range("D14").Value = "'5"
$ws.Range("E14").Value = 20
$ws.Range("L14").Value = "'8"
$ws.Range("AB14").Value = "'10"
$ws.Range("D15").Value = "'2"
$ws.Range("L15").Value = "'8"
$ws.Range("AB15").Value = "'9"
$ws.Range("AC15").Value = 0.59

$ws.Range("E4").Select()
